$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "<head>Preventing <bp>teats</bp>..." -> "<head>To prevent <bp>teats</bp>..."
#   Original runs: "P" | "revent" | "ing "   (all w:rtl="0"; middle one also has color="000000")
#   Target runs:   "To prevent" | " "        (both just w:rtl="0", no color)
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("Preventing ")
if ($idx -lt 0) { throw "could not locate 'Preventing ' in the document" }

# 1a) Turn "P" + "revent" into a single run "To prevent" (keeps the first run's
#     formatting, i.e. rtl=0 / no color, same as the target).
$rngPrevent = $d.Range($idx, $idx + "Prevent".Length)
$rngPrevent.Find.Execute("Prevent", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "To prevent", 2)

# 1b) Drop a bookmark right at the new run boundary so that the engine's
#     adjacent-run coalescing (which would otherwise merge "To prevent" with
#     the following "ing " run, since both end up formatted identically)
#     cannot fire while we edit the trailing run.
$text = $d.Content.Text
$idx = $text.IndexOf("To preventing")
$boundary = $d.Range($idx + "To prevent".Length, $idx + "To prevent".Length)
$d.Bookmarks.Add("zzEditBoundary", $boundary) | Out-Null

# 1c) Shrink "ing " down to just " " (its own run, untouched formatting).
$text = $d.Content.Text
$idx = $text.IndexOf("To preventing")
$rngIng = $d.Range($idx + "To prevent".Length, $idx + "To preventing ".Length)
$rngIng.Find.Execute("ing ", $true, $false, $false, $false, $false, `
                      $true, 1, $false, " ", 2)

# 1d) Remove the helper bookmark; the two runs stay distinct.
if ($d.Bookmarks.Exists("zzEditBoundary")) {
    $d.Bookmarks.Item("zzEditBoundary").Delete()
}

# ---------------------------------------------------------------------
# Edit 2 & 3: "<head>Against the bruising of <bp>eyes</bp></head>"
#           -> "<head>Against bruising of the <bp>eyes</bp></head>"
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("Against the bruising of ")
if ($idx -lt 0) { throw "could not locate 'Against the bruising of ' in the document" }
$rngAgainst = $d.Range($idx, $idx + "Against the bruising of ".Length)
$rngAgainst.Find.Execute("Against the ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Against ", 2)

$text = $d.Content.Text
$idx = $text.IndexOf("Against bruising of ")
if ($idx -lt 0) { throw "could not locate 'Against bruising of ' in the document" }
$rngOf = $d.Range($idx, $idx + "Against bruising of ".Length)
$rngOf.Find.Execute("of ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "of the ", 2)
